$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the participation value for the Quentin Lintz / Implementation row (row 13)
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = "Implementation"

# Move the active selection to G14, matching the saved cursor position
$ws.Range("G14").Select()
